# Refresh the "number of people interested" counts (column F) on every sheet
# so the workbook matches the newly generated gh-pages data snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 257
$ws.Range("F3").Value = 2742
$ws.Range("F5").Value = 942
$ws.Range("F6").Value = 37
$ws.Range("F7").Value = 2377
$ws.Range("F8").Value = 1853
$ws.Range("F9").Value = 220
$ws.Range("F11").Value = 2505
$ws.Range("F12").Value = 561
$ws.Range("F13").Value = 247
$ws.Range("F17").Value = 121
$ws.Range("F18").Value = 9344
$ws.Range("F20").Value = 7254
$ws.Range("F21").Value = 11820
$ws.Range("F26").Value = 566
$ws.Range("F27").Value = 2637
$ws.Range("F28").Value = 236
$ws.Range("F29").Value = 202
$ws.Range("F30").Value = 2592
$ws.Range("F31").Value = 777
$ws.Range("F33").Value = 4533
$ws.Range("F34").Value = 971
$ws.Range("F36").Value = 48
$ws.Range("F37").Value = 541

$ws = $wb.Worksheets.Item(2)
$ws.Range("F17").Value = 67
$ws.Range("F22").Value = 17
$ws.Range("F24").Value = 14
$ws.Range("F27").Value = 7

$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 166

$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 257
$ws.Range("F5").Value = 2742
$ws.Range("F8").Value = 942
$ws.Range("F9").Value = 37
$ws.Range("F11").Value = 2377
$ws.Range("F13").Value = 1853
$ws.Range("F14").Value = 220
$ws.Range("F15").Value = 2505
$ws.Range("F17").Value = 561
$ws.Range("F18").Value = 247
$ws.Range("F22").Value = 121
$ws.Range("F23").Value = 9344
$ws.Range("F25").Value = 7254
$ws.Range("F26").Value = 11820
$ws.Range("F32").Value = 566
$ws.Range("F34").Value = 2637
$ws.Range("F37").Value = 236
$ws.Range("F38").Value = 202
$ws.Range("F40").Value = 4533
$ws.Range("F41").Value = 67
$ws.Range("F45").Value = 541
$ws.Range("F46").Value = 14
$ws.Range("F49").Value = 7
